$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.864.28'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '1.667.08'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.46'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.532'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.47%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("E9").Value = '  +0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.18'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.24%  '
$ws.Range("E11").Value = '  +3.75%  '
$ws.Range("D12").Value = '1.903.68'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '1.675.13'
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.523'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.97'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("D17").Value = '26.876.39'
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.96'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -4.03%  '
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D20").Value = '0.0₃0733'
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.46'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("E24").Value = '  -1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.57'
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.12'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.88'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -0.45%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.48%  '
$ws.Range("D33").Value = '1.467.08'
$ws.Range("E33").Value = '  -3.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.55%  '
$ws.Range("E35").Value = '  +2.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.41'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("E37").Value = '  -0.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.897'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.55%  '
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.84'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("E42").Value = '  -2.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.979'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +6.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.83'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.81%  '
$ws.Range("D45").Value = '1.812.98'
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.777'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.21'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("E49").Value = '  +2.91%  '
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("E51").Value = '  +0.61%  '
